# Fix mojibake "Â±" (chars 194,177) -> "±" (char 177) in columns B, C, D
# for rows 2 through 17 on the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$capA = [char]194
$plusMinus = [char]177
$mojibake = "$capA$plusMinus"
$replacement = "$plusMinus"

for ($row = 2; $row -le 17; $row++) {
    foreach ($col in @("B", "C", "D")) {
        $cell = $ws.Range("$col$row")
        $value = $cell.Value2
        if ($null -ne $value -and $value -is [string] -and $value.Contains($mojibake)) {
            $cell.Value = $value.Replace($mojibake, $replacement)
        }
    }
}
